$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 'Disciplinar'
$ws.Range("G4").Value = 1

$ws.Range("F7").Value = 'Fund. Obligatoria'
$ws.Range("G7").Value = 1

$ws.Range("F8").Value = 'Fund. Obligatoria'
$ws.Range("G8").Value = 2

$ws.Range("F9").Value = 'Fund. Obligatoria'
$ws.Range("G9").Value = 2

$ws.Range("F11").Value = 'Fund. Obligatoria'
$ws.Range("G11").Value = 2

$ws.Range("F14").Value = 'Fund. Obligatoria'
$ws.Range("G14").Value = 3

$ws.Range("F15").Value = 'Fund. Obligatoria'
$ws.Range("G15").Value = 4

$ws.Range("F22").Value = 'Fund. Obligatoria'
$ws.Range("G22").Value = 2

$ws.Range("F23").Value = 'Fund. Obligatoria'
$ws.Range("G23").Value = 4

$ws.Range("F65").Value = 'Nivelación'
$ws.Range("G65").Value = 1

$ws.Range("F66").Value = 'Fund. Obligatoria'
$ws.Range("G66").Value = 1

$ws.Range("F67").Value = 'Nivelación'
$ws.Range("G67").Value = 1

$ws.Range("F68").Value = 'Fund. Obligatoria'
$ws.Range("G68").Value = 1

$ws.Range("F70").Value = 'Disciplinar'
$ws.Range("G70").Value = 1

$ws.Range("F71").Value = 'Fund. Obligatoria'
$ws.Range("G71").Value = 2

$ws.Range("F72").Value = 'Fund. Obligatoria'
$ws.Range("G72").Value = 1

$ws.Range("F73").Value = 'Fund. Obligatoria'
$ws.Range("G73").Value = 4

$ws.Range("F74").Value = 'Fund. Obligatoria'
$ws.Range("G74").Value = 2

$ws.Range("F76").Value = 'Disciplinar'
$ws.Range("G76").Value = 5

$ws.Range("F77").Value = 'Fund. Obligatoria'
$ws.Range("G77").Value = 3

$ws.Range("F78").Value = 'Fund. Obligatoria'
$ws.Range("G78").Value = 2

$ws.Range("F79").Value = 'Fund. Obligatoria'
$ws.Range("G79").Value = 2

$ws.Range("F80").Value = 'Disciplinar'
$ws.Range("G80").Value = 4

$ws.Range("F82").Value = 'Fund. Obligatoria'
$ws.Range("G82").Value = 3

$ws.Range("F83").Value = 'Fund. Obligatoria'
$ws.Range("G83").Value = 2

$ws.Range("F84").Value = 'Fund. Obligatoria'
$ws.Range("G84").Value = 3

$ws.Range("F85").Value = 'Fund. Obligatoria'
$ws.Range("G85").Value = 4

$ws.Range("F86").Value = 'Disciplinar'
$ws.Range("G86").Value = 3

$ws.Range("F87").Value = 'Fund. Obligatoria'
$ws.Range("G87").Value = 4

$ws.Range("F88").Value = 'Fund. Obligatoria'
$ws.Range("G88").Value = 3

$ws.Range("F90").Value = 'Fund. Obligatoria'
$ws.Range("G90").Value = 5

$ws.Range("F91").Value = 'Fund. Obligatoria'
$ws.Range("G91").Value = 4

$ws.Range("F92").Value = 'Fund. Obligatoria'
$ws.Range("G92").Value = 2

$ws.Range("F93").Value = 'Disciplinar'
$ws.Range("G93").Value = 6

$ws.Range("F95").Value = 'Disciplinar'
$ws.Range("G95").Value = 5

$ws.Range("F97").Value = 'Disciplinar'
$ws.Range("G97").Value = 4

$ws.Range("F98").Value = 'Disciplinar'
$ws.Range("G98").Value = 6

$ws.Range("F99").Value = 'Disciplinar'
$ws.Range("G99").Value = 7

$ws.Range("F100").Value = 'Disciplinar'
$ws.Range("G100").Value = 6

$ws.Range("F101").Value = 'Disciplinar'
$ws.Range("G101").Value = 5

$ws.Range("F102").Value = 'Disciplinar'
$ws.Range("G102").Value = 6

$ws.Range("F103").Value = 'Disciplinar'
$ws.Range("G103").Value = 6

$ws.Range("F107").Value = 'Disciplinar'
$ws.Range("G107").Value = 7

$ws.Range("F108").Value = 'Disciplinar'
$ws.Range("G108").Value = 7

$ws.Range("F109").Value = 'Disciplinar'
$ws.Range("G109").Value = 7

$ws.Range("F111").Value = 'Disciplinar'
$ws.Range("G111").Value = 7

$ws.Range("F112").Value = 'Disciplinar'
$ws.Range("G112").Value = 7

$ws.Range("F113").Value = 'Disciplinar'
$ws.Range("G113").Value = 8

$ws.Range("F116").Value = 'Disciplinar'
$ws.Range("G116").Value = 8

$ws.Range("F117").Value = 'Disciplinar'
$ws.Range("G117").Value = 5

$ws.Range("F118").Value = 'Disciplinar'
$ws.Range("G118").Value = 8

$ws.Range("F119").Value = 'Disciplinar'
$ws.Range("G119").Value = 9

$ws.Range("F121").Value = 'Optativa de Producción'
$ws.Range("G121").Value = 9

$ws.Range("F123").Value = 'Optativa de Producción'

$ws.Range("F131").Value = 'Optativa de Producción'
